# Bug fix: the model's start year was incorrectly set to 2015; correct it to 2018.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2018

# Leave the selection where the user would naturally land after editing the
# StartYear cell (matches the saved sheetView selection in the fixed file).
[void]$ws.Range("B3").Select()
